$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as text (e.g. "30.549.94", "0.4830")
# in the original workbook. Excel auto-converts numeric-looking strings to
# real numbers on assignment, which would lose formatting such as trailing
# zeros or thousands separators. To avoid that, each Price cell is temporarily
# switched to a text number format before its value is written, then restored
# to the default (Normal) style so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.549.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.917.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4830"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2892"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06698"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "111.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.90%  "
$ws.Range("E11").Value = "  +2.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.920.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07579"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.276"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6678"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "295.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.529.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007595"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.550"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.164.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.445"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.459"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("E26").Value = "  -2.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.090"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1070"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.443"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.132"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.050"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05009"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7402"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.136"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.720"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02023"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.683"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.013"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4419"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8666"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.826"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.17%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.217"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.223"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.76%  "
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2518"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.15%  "
